# Update the Leave Card ("Sheet1") with new leave-credit entries for 2023.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 19: VL(1-0-0) entry for 2022-11-30, 1.25 days earned ---
$ws.Range("C19").Value = 1.25

# --- Row 20: new dated entry (2023-12-31), 1.25 days earned ---
$ws.Range("A20").Value = 45291
$ws.Range("C20").Value = 1.25

# --- Row 21: "2023" year header row (mirrors the "2022" header at A10) ---
$ws.Range("A21").Value = "2023"
$ws.Range("A21").Style = $ws.Range("A10").Style

# --- Rows 22-26: monthly VL(1-0-0) credits for 2023, 1.25 days each ---
$ws.Range("A22").Value = 44957
$ws.Range("C22").Value = 1.25

$ws.Range("A23").Value = 44985
$ws.Range("C23").Value = 1.25

$ws.Range("A24").Value = 45016
$ws.Range("C24").Value = 1.25

$ws.Range("A25").Value = 45046
$ws.Range("C25").Value = 1.25

$ws.Range("A26").Value = 45077
$ws.Range("C26").Value = 1.25

# --- Row 27: 2023-06-30 entry, "VL(1-0-0)" particulars, 1 day absence, remarks date ---
$ws.Range("A27").Value = 45107
$ws.Range("B27").Value = "VL(1-0-0)"
$ws.Range("D27").Value = 1
$ws.Range("K27").Value = 45108
$ws.Range("K27").Style = $ws.Range("K18").Style

# --- Rows 28-33: remaining month-end dates for 2023 ---
$ws.Range("A28").Value = 45138
$ws.Range("A29").Value = 45169
$ws.Range("A30").Value = 45199
$ws.Range("A31").Value = 45230
$ws.Range("A32").Value = 45231
$ws.Range("A33").Value = 45232
